# The workbook's "in" sheet lists model results grouped by offense
# category (General, Violent, Drug, Property, Felony, Misdemeanor),
# with a blank row left between each group. Previously there was no
# blank row between "General" (row 3) and "Violent" (row 4). This
# edit inserts a blank row after row 3 so the layout is consistent
# with the rest of the table, which pushes the "Violent", "Drug",
# "Property", "Felony" and "Misdemeanor" blocks (and their data) down
# by one row each (4->5, 6->7, 8->9, 10->11, 12->13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4").EntireRow.Insert() | Out-Null

# Reflect the last active selection recorded in the saved workbook.
$ws.Range("D10").Select() | Out-Null
